$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (style) from N1 into the two new header cells
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$oValues = @(
14,-0.2114643095450728,0.2806349495657003,0.2473511196078822,-0.03033007942995938,0.2293198790880224,-0.3700356111966396,-0.1413961698196282,-0.2819806693959612,0.4300107267053127,-0.176727984093035,-0.007243021530742952,-0.02024715154380544,0.20396057915699,0.1084671529106905,0.5579117559550353,0.6362965422861516,-0.06815352160258968,0.4115014975651828,0.4588847377591326,0.6728430254908768,0.4425085211407857,-0.0113000032422991,2.129520835592704,0.2904976299505371,0.1506317619940473,0.05399768676252616,0.7954418179133347,1.940548621359262,0.6386472370929314,-0.4983202327706122,0.5284334960304413,0.7367409381924135,-0.8579614069301367,0.7745650181102949,0.7133374342039093,0.6762752632840786,0.6245836616786016,0.5770265864720949,0.7392385292386937,0.5499278257552696,0.5708037534225769,0.6505680392745654,0.6650497093602901,0.622719689574091,-1.288991576725659,-1.006854996677864,-0.880998849997614,-0.6427838223261589,-0.06159562467068196,-0.8684286529132981,-0.8684286529132981,-1.142703383246583,-0.175790836190442,-1.025965455424598,-0.894838901547386,-0.9028397994521938,-1.047861740587387,-0.7800699864657531,-0.4030377529607465,0.3877615184819895,-1.153323077767894,-0.5616580007143084,-0.8368302087400017,-0.02155160459750484,-0.7314242858321528,-0.7003664362236856
)

$pValues = @(
15,-0.2105307011609714,0.2804782342884273,0.2495978380119809,-0.03184136088923693,0.2286045838077772,-0.368788297950749,-0.1369774069469465,-0.2788143977689841,0.4288048919811859,-0.1763990854416434,-0.01241214794415214,-0.02248655969775492,0.1954692355642895,0.09458394629827144,0.5433729377062484,0.6196207302970599,-0.07636216074010468,0.4027274116305298,0.4420835603505659,0.6594159441139804,0.4281838752811294,-0.02378779627998907,1.933480563269069,0.2884788020759164,0.1423364354741517,0.04474693881610786,0.7878710865875564,1.808452676035163,0.6317236554716106,-0.4982459878878621,0.52216684504216,0.7339635203510406,-0.8615942352856631,0.7764889788772975,0.7177841561524801,0.6807298854737832,0.6248803826165612,0.5794582359511189,0.741245263783091,0.5529337956407395,0.5736166647592651,0.6521617104953988,0.6686488913861264,0.6300826381578328,-1.291101303561887,-1.008407106267946,-0.8804238944387288,-0.6425418421393924,-0.06295277074671651,-0.8672392379887609,-0.8672392379887609,-1.14308967652738,-0.1758729849932592,-1.02785836688306,-0.8978345750150249,-0.9084518483657634,-1.05259704072894,-0.7805983542856583,-0.4047295025405627,0.3876488883810655,-1.159758715534909,-0.5577671524865098,-0.8369824073405764,-0.02267873477588057,-0.7379955461768175,-0.7094190390418398
)

for ($i = 0; $i -lt 67; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 15).Value = [double]$oValues[$i]
    $ws.Cells.Item($row, 16).Value = [double]$pValues[$i]
}
